$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "img" column (E) by shifting "images" (F) and "desc" (G)
#    one column to the left, then deleting the now-duplicated last column.
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 3; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value()
    $gVal = $ws.Cells.Item($r, 7).Value()
    $ws.Cells.Item($r, 5).Value = $fVal
    $ws.Cells.Item($r, 6).Value = $gVal
}
$ws.Columns.Item(7).Delete()

# ---------------------------------------------------------------------------
# 2. Update row 3 (id = 1) with the new article content.
# ---------------------------------------------------------------------------
$titleText = "SGRU and CEPT to Host ""Smart Building Data Analytics"" Workshop on AI-Driven Energy Management"
$dateText = "12-13 Mar 2020"
$categoryText = "Activity"
$linksText = "https://drive.google.com/file/d/1v0jfCB_yGTld_LgZBUqY-COEVeQh0lxx/view?usp=drive_link, https://drive.google.com/file/d/1tBnLfQddPVHBfKwgbfxOK96UV8a-8ybl/view?usp=drive_link, https://drive.google.com/file/d/1I_AwRAwLcxNczYM8TVcGQ3qCZ17_yPQB/view?usp=drive_link, https://drive.google.com/file/d/1GiMyGgwBBY_c4trGpkKRJq6WyVVCRk4g/view?usp=drive_link, https://drive.google.com/file/d/1tBGB_ecQerpYyf5zEja830-NpDnlvZQE/view?usp=drive_link"
$descText = "Building Energy Management Systems (BEMS) have evolved significantly through the integration of IoT and AI, enabling sophisticated energy analysis and optimization that is now accessible to both large organizations and individual users. To address these advancements, the Smart Grid Research Unit (SGRU) and the Center of Excellence in Electrical Power Technology (CEPT) at Chulalongkorn University are hosting the ""Smart Building Data Analytics"" course on March 12-13, 2020. This program is designed to equip participants with practical expertise in applying AI and Big Data analytics to enhance energy management, aligning with current global technological trends."

$ws.Range("B3").Value = $titleText
$ws.Range("C3").Value = $dateText
$ws.Range("D3").Value = $categoryText
$ws.Range("E3").Value = $linksText
$ws.Range("F3").Value = $descText

# ---------------------------------------------------------------------------
# 3. Formatting - header row wraps text, body rows get a wrapping style, and
#    the new "images" cell on row 3 (Google Drive links) gets its own font.
# ---------------------------------------------------------------------------
$ws.Range("A1:F1").WrapText = $true

$ws.Range("A2:F3").WrapText = $true

$f = $ws.Range("E3").Font
$f.Size = 9
$f.Name = "Segoe UI"
$f.Color = 10058872

# ---------------------------------------------------------------------------
# 4. Row heights for the two data rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 150

# ---------------------------------------------------------------------------
# 5. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.3
$ws.Columns.Item(2).ColumnWidth = 35.17
$ws.Columns.Item(3).ColumnWidth = 11.3
$ws.Columns.Item(4).ColumnWidth = 17.74
$ws.Columns.Item(5).ColumnWidth = 53.88
$ws.Columns.Item(6).ColumnWidth = 71.02

# ---------------------------------------------------------------------------
# 6. Freeze the header row and set the view/selection state.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E8").Select()
